$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.566.37"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.831.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.03%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "429.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.823.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.45%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.727"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.04%  "

# Row 11
$ws.Range("E11").Value = "  -9.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000365"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.82%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.439.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.33%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +17.15%  "

# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.39%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.831.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.83%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.95%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.008.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21
$ws.Range("E21").Value = "  -6.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "409.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.99%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.08%  "

# Row 25
$ws.Range("E25").Value = "  -4.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.00%  "

# Row 27
$ws.Range("E27").Value = "  +12.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.88%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "687.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.86%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.68%  "

# Row 32
$ws.Range("E32").Value = "  -2.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.47%  "

# Row 34
$ws.Range("E34").Value = "  -1.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.151"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.62%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.98%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0791"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.85%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.77%  "

# Row 41
$ws.Range("E41").Value = "  -8.16%  "

# Row 42
$ws.Range("E42").Value = "  +0.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.81%  "

# Row 44
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.68%  "

# Row 45
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.07%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.15%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.16%  "

# Row 51
$ws.Range("E51").Value = "  -5.02%  "
